$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "64.972.58"
$ws.Cells.Item(2, 5).Value = "  +0.93%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.366.38"
$ws.Cells.Item(3, 5).Value = "  +0.83%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'554.31"
$ws.Cells.Item(5, 5).Value = "  +0.42%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'173.86"
$ws.Cells.Item(6, 5).Value = "  -0.66%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +2.31%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.359.45"
$ws.Cells.Item(8, 5).Value = "  +0.83%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.01%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +7.74%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.55%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'53.59"
$ws.Cells.Item(12, 5).Value = "  -1.90%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +3.32%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'9.12"
$ws.Cells.Item(14, 5).Value = "  +0.96%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.902.68"
$ws.Cells.Item(15, 5).Value = "  +0.71%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +2.36%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.47%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.362.65"
$ws.Cells.Item(18, 5).Value = "  +0.76%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "WrappedBTC"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(19, 4).Value = "64.964.82"
$ws.Cells.Item(19, 5).Value = "  +1.09%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Value = "'11.87"
$ws.Cells.Item(20, 5).Value = "  +1.20%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'0.996"
$ws.Cells.Item(21, 5).Value = "  +1.80%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'449.95"
$ws.Cells.Item(22, 5).Value = "  +4.00%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'4.94"
$ws.Cells.Item(23, 5).Value = "  -2.85%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.05%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'87.07"
$ws.Cells.Item(25, 5).Value = "  +3.30%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'13.65"
$ws.Cells.Item(26, 5).Value = "  +2.08%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'10.72"
$ws.Cells.Item(27, 5).Value = "  +0.06%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +1.21%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'8.66"
$ws.Cells.Item(29, 5).Value = "  -0.74%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'31.15"
$ws.Cells.Item(30, 5).Value = "  +4.84%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -1.92%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'62.90"
$ws.Cells.Item(32, 5).Value = "  +8.04%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'11.43"
$ws.Cells.Item(33, 5).Value = "  -0.39%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'575.74"
$ws.Cells.Item(34, 5).Value = "  -0.61%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.20%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.03%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'3.63"
$ws.Cells.Item(37, 5).Value = "  +4.05%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.14%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'35.61"
$ws.Cells.Item(39, 5).Value = "  +0.13%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +1.06%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "0.0₃0740"
$ws.Cells.Item(41, 5).Value = "  -1.32%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "3.079.29"
$ws.Cells.Item(42, 5).Value = "  -0.80%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +1.95%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -1.38%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Fetch.AI"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(46, 4).Value = "'2.43"
$ws.Cells.Item(46, 5).Value = "  -1.02%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "ApeXProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(47, 4).Value = "'3.15"
$ws.Cells.Item(47, 5).Value = "  -2.05%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(48, 4).Value = "'0.999"
$ws.Cells.Item(48, 5).Value = "  -0.05%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Monero"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(49, 4).Value = "'141.49"
$ws.Cells.Item(49, 5).Value = "  +3.98%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'2.54"
$ws.Cells.Item(50, 5).Value = "  -2.15%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'8.27"
$ws.Cells.Item(51, 5).Value = "  -0.06%  "
